$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply column-A formatting (bold, centered, thin border) to the newly added rows
# by copying the format from an existing styled A-column cell (A2), so we reuse
# the existing style index instead of creating new ones.
$ws.Range("A2").Copy()
$ws.Range("A59:A81").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 54
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 3.203383214053351
$ws.Range("C54").Value = 4923.6
$ws.Range("D54").Value = 0.01535458685751464
$ws.Range("E54").Value = 23.6
$ws.Range("F54").Value = 126
$ws.Range("G54").Value = 'Châllénger '
$ws.Range("H54").Value = 'SOLO'
$ws.Range("I54").Value = 0.08197787898503578
$ws.Range("J54").Value = 3.4
$ws.Range("K54").Value = 0.0022121014964216

# Row 55
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 1.775319622012229
$ws.Range("C55").Value = 3193.8
$ws.Range("D55").Value = 0.01634241245136187
$ws.Range("E55").Value = 29.4
$ws.Range("F55").Value = 41
$ws.Range("G55").Value = 'Portgas D Åce '
$ws.Range("H55").Value = 'SOLO'
$ws.Range("I55").Value = 0.02279043913285158
$ws.Range("J55").Value = 4.4
$ws.Range("K55").Value = 0.002445803224013341

# Row 56
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 4.58653314427508
$ws.Range("C56").Value = 8017.8
$ws.Range("D56").Value = 0.03247285569866215
$ws.Range("E56").Value = 56.2
$ws.Range("F56").Value = 187.4
$ws.Range("G56").Value = 'LS DUFFY'
$ws.Range("H56").Value = 'SOLO'
$ws.Range("I56").Value = 0.1066721808657292
$ws.Range("J56").Value = 5.4
$ws.Range("K56").Value = 0.003094384707287933

# Row 57
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 3.650427670597207
$ws.Range("C57").Value = 5630
$ws.Range("D57").Value = 0.02674180615384204
$ws.Range("E57").Value = 44.4
$ws.Range("F57").Value = 227.4
$ws.Range("G57").Value = 'BigFather Rengar'
$ws.Range("H57").Value = 'SOLO'
$ws.Range("I57").Value = 0.1559623334516667
$ws.Range("J57").Value = 11.4
$ws.Range("K57").Value = 0.00728984276207101

# Row 58
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 1.657101658255227
$ws.Range("C58").Value = 2298.4
$ws.Range("D58").Value = 0.02033165104542177
$ws.Range("E58").Value = 28.2
$ws.Range("F58").Value = 69.2
$ws.Range("G58").Value = '19 fotsiny adc'
$ws.Range("H58").Value = 'NONE'
$ws.Range("I58").Value = 0.04989185291997116
$ws.Range("J58").Value = 0.2
$ws.Range("K58").Value = 0.0001441961067051189

# Row 59
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 2.219149527515286
$ws.Range("C59").Value = 3992.25
$ws.Range("D59").Value = 0.02042801556420233
$ws.Range("E59").Value = 36.75
$ws.Range("F59").Value = 51.25
$ws.Range("G59").Value = 'Portgas D Åce '
$ws.Range("H59").Value = 'SOLO'
$ws.Range("I59").Value = 0.02848804891606448
$ws.Range("J59").Value = 5.5
$ws.Range("K59").Value = 0.003057254030016676

# Row 60
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 5.73316643034385
$ws.Range("C60").Value = 10022.25
$ws.Range("D60").Value = 0.04059106962332769
$ws.Range("E60").Value = 70.25
$ws.Range("F60").Value = 234.25
$ws.Range("G60").Value = 'LS DUFFY'
$ws.Range("H60").Value = 'SOLO'
$ws.Range("I60").Value = 0.1333402260821616
$ws.Range("J60").Value = 6.75
$ws.Range("K60").Value = 0.003867980884109916

# Row 61
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 2.071377072819034
$ws.Range("C61").Value = 2873
$ws.Range("D61").Value = 0.02541456380677722
$ws.Range("E61").Value = 35.25
$ws.Range("F61").Value = 86.5
$ws.Range("G61").Value = '19 fotsiny adc'
$ws.Range("H61").Value = 'NONE'
$ws.Range("I61").Value = 0.06236481614996395
$ws.Range("J61").Value = 0.25
$ws.Range("K61").Value = 0.0001802451333813987

# Row 62
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 6.451745292628223
$ws.Range("C62").Value = 11335.4
$ws.Range("D62").Value = 0.03638221087160513
$ws.Range("E62").Value = 65.8
$ws.Range("F62").Value = 293
$ws.Range("G62").Value = 'BigFather Rengar'
$ws.Range("H62").Value = 'SOLO'
$ws.Range("I62").Value = 0.1866635435052542
$ws.Range("J62").Value = 13.6
$ws.Range("K62").Value = 0.008109015746182139

# Row 63
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = 4.768732242866904
$ws.Range("C63").Value = 7445
$ws.Range("D63").Value = 0.03195378050521529
$ws.Range("E63").Value = 50.16666666666666
$ws.Range("F63").Value = 190.1666666666667
$ws.Range("G63").Value = 'Châllénger '
$ws.Range("H63").Value = 'SOLO'
$ws.Range("I63").Value = 0.1218115490704444
$ws.Range("J63").Value = 4.833333333333333
$ws.Range("K63").Value = 0.003099699320719843

# Row 64
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = 1.479433018343524
$ws.Range("C64").Value = 2661.5
$ws.Range("D64").Value = 0.01361867704280156
$ws.Range("E64").Value = 24.5
$ws.Range("F64").Value = 34.16666666666666
$ws.Range("G64").Value = 'Portgas D Åce '
$ws.Range("H64").Value = 'SOLO'
$ws.Range("I64").Value = 0.01899203261070965
$ws.Range("J64").Value = 3.666666666666667
$ws.Range("K64").Value = 0.002038169353344451

# Row 65
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 3.822110953562567
$ws.Range("C65").Value = 6681.5
$ws.Range("D65").Value = 0.02706071308221846
$ws.Range("E65").Value = 46.83333333333334
$ws.Range("F65").Value = 156.1666666666667
$ws.Range("G65").Value = 'LS DUFFY'
$ws.Range("H65").Value = 'SOLO'
$ws.Range("I65").Value = 0.08889348405477437
$ws.Range("J65").Value = 4.5
$ws.Range("K65").Value = 0.002578653922739944

# Row 66
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 7.109574936737488
$ws.Range("C66").Value = 12432.33333333333
$ws.Range("D66").Value = 0.0492776887849563
$ws.Range("E66").Value = 87.5
$ws.Range("F66").Value = 326.8333333333333
$ws.Range("G66").Value = 'BigFather Rengar'
$ws.Range("H66").Value = 'SOLO'
$ws.Range("I66").Value = 0.2035312852870734
$ws.Range("J66").Value = 17.83333333333333
$ws.Range("K66").Value = 0.01053000296511503

# Row 67
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 5.251755500694689
$ws.Range("C67").Value = 9529.166666666666
$ws.Range("D67").Value = 0.06558663539837613
$ws.Range("E67").Value = 113.5
$ws.Range("F67").Value = 266.6666666666667
$ws.Range("G67").Value = 'JaIisco'
$ws.Range("H67").Value = 'SOLO'
$ws.Range("I67").Value = 0.1534088627803535
$ws.Range("J67").Value = 8.333333333333334
$ws.Range("K67").Value = 0.004630813408796165

# Row 68
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 3.393041794365421
$ws.Range("C68").Value = 4321.833333333333
$ws.Range("D68").Value = 0.03561639816772885
$ws.Range("E68").Value = 45.83333333333334
$ws.Range("F68").Value = 98
$ws.Range("G68").Value = '19 fotsiny adc'
$ws.Range("H68").Value = 'NONE'
$ws.Range("I68").Value = 0.0753000669539336
$ws.Range("J68").Value = 1.5
$ws.Range("K68").Value = 0.001234990624037989

# Row 69
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 1.548397645519948
$ws.Range("C69").Value = 2367.5
$ws.Range("D69").Value = 0.02043819489862655
$ws.Range("E69").Value = 31.25
$ws.Range("F69").Value = 97.5
$ws.Range("G69").Value = 'Booogeyman'
$ws.Range("H69").Value = 'DUO_CARRY'
$ws.Range("I69").Value = 0.06376716808371484
$ws.Range("J69").Value = 4
$ws.Range("K69").Value = 0.002616088947024199

# Row 70
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = 1.158510638297872
$ws.Range("C70").Value = 1089
$ws.Range("D70").Value = 0.02154255319148936
$ws.Range("E70").Value = 20.25
$ws.Range("F70").Value = 107
$ws.Range("G70").Value = 'Cantare'
$ws.Range("H70").Value = 'DUO_CARRY'
$ws.Range("I70").Value = 0.1138297872340426
$ws.Range("J70").Value = 0.75
$ws.Range("K70").Value = 0.0007978723404255319

# Row 71
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = 3.703620136892431
$ws.Range("C71").Value = 6127.25
$ws.Range("D71").Value = 0.04150350867466872
$ws.Range("E71").Value = 67
$ws.Range("F71").Value = 168.5
$ws.Range("G71").Value = 'Poppy Gods'
$ws.Range("H71").Value = 'DUO_CARRY'
$ws.Range("I71").Value = 0.1027973605081667
$ws.Range("J71").Value = 8.5
$ws.Range("K71").Value = 0.005183907954944446

# Row 72
$ws.Range("A72").Value = 70
$ws.Range("B72").Value = 1.933228192857932
$ws.Range("C72").Value = 2893.4
$ws.Range("D72").Value = 0.02065910352140298
$ws.Range("E72").Value = 31.2
$ws.Range("F72").Value = 184.8
$ws.Range("G72").Value = 'Booogeyman'
$ws.Range("H72").Value = 'DUO_CARRY'
$ws.Range("I72").Value = 0.1252319415552276
$ws.Range("J72").Value = 9.199999999999999
$ws.Range("K72").Value = 0.006262433353588783

# Row 73
$ws.Range("A73").Value = 71
$ws.Range("B73").Value = 2.278052291732825
$ws.Range("C73").Value = 1142.8
$ws.Range("D73").Value = 0.0271842913094104
$ws.Range("E73").Value = 18.2
$ws.Range("F73").Value = 122.4
$ws.Range("G73").Value = 'Cantare'
$ws.Range("H73").Value = 'DUO_CARRY'
$ws.Range("I73").Value = 0.2741484069016619
$ws.Range("J73").Value = 0.6
$ws.Range("K73").Value = 0.0006382978723404255

# Row 74
$ws.Range("A74").Value = 72
$ws.Range("B74").Value = 2.962896109513945
$ws.Range("C74").Value = 4901.8
$ws.Range("D74").Value = 0.03320280693973497
$ws.Range("E74").Value = 53.6
$ws.Range("F74").Value = 134.8
$ws.Range("G74").Value = 'Poppy Gods'
$ws.Range("H74").Value = 'DUO_CARRY'
$ws.Range("I74").Value = 0.08223788840653333
$ws.Range("J74").Value = 6.8
$ws.Range("K74").Value = 0.004147126363955557

# Row 75
$ws.Range("A75").Value = 73
$ws.Range("B75").Value = 6.710934922957702
$ws.Range("C75").Value = 11434.2
$ws.Range("D75").Value = 0.0404788960127518
$ws.Range("E75").Value = 62.2
$ws.Range("F75").Value = 338.2
$ws.Range("G75").Value = 'Ithryn L'
$ws.Range("H75").Value = 'DUO_CARRY'
$ws.Range("I75").Value = 0.2035710766994756
$ws.Range("J75").Value = 8.199999999999999
$ws.Range("K75").Value = 0.004878862988729487

# Row 76
$ws.Range("A76").Value = 74
$ws.Range("B76").Value = 8.922961909838907
$ws.Range("C76").Value = 15499.4
$ws.Range("D76").Value = 0.0439786500407446
$ws.Range("E76").Value = 78.59999999999999
$ws.Range("F76").Value = 418.4
$ws.Range("G76").Value = 'BigFather Rengar'
$ws.Range("H76").Value = 'SOLO'
$ws.Range("I76").Value = 0.2610849084904173
$ws.Range("J76").Value = 18.2
$ws.Range("K76").Value = 0.01083898607259163

# Row 77
$ws.Range("A77").Value = 75
$ws.Range("B77").Value = 7.352278870698461
$ws.Range("C77").Value = 13452.8
$ws.Range("D77").Value = 0.03998785955724499
$ws.Range("E77").Value = 73.40000000000001
$ws.Range("F77").Value = 290.2
$ws.Range("G77").Value = 'BigFather Rengar'
$ws.Range("H77").Value = 'SOLO'
$ws.Range("I77").Value = 0.1626965738779845
$ws.Range("J77").Value = 14.6
$ws.Range("K77").Value = 0.008076131122476515

# Row 78
$ws.Range("A78").Value = 80
$ws.Range("B78").Value = 3.203383214053351
$ws.Range("C78").Value = 4923.6
$ws.Range("D78").Value = 0.01535458685751464
$ws.Range("E78").Value = 23.6
$ws.Range("F78").Value = 126
$ws.Range("G78").Value = 'Châllénger '
$ws.Range("H78").Value = 'SOLO'
$ws.Range("I78").Value = 0.08197787898503578
$ws.Range("J78").Value = 3.4
$ws.Range("K78").Value = 0.0022121014964216

# Row 79
$ws.Range("A79").Value = 81
$ws.Range("B79").Value = 1.775319622012229
$ws.Range("C79").Value = 3193.8
$ws.Range("D79").Value = 0.01634241245136187
$ws.Range("E79").Value = 29.4
$ws.Range("F79").Value = 41
$ws.Range("G79").Value = 'Portgas D Åce '
$ws.Range("H79").Value = 'SOLO'
$ws.Range("I79").Value = 0.02279043913285158
$ws.Range("J79").Value = 4.4
$ws.Range("K79").Value = 0.002445803224013341

# Row 80
$ws.Range("A80").Value = 82
$ws.Range("B80").Value = 4.58653314427508
$ws.Range("C80").Value = 8017.8
$ws.Range("D80").Value = 0.03247285569866215
$ws.Range("E80").Value = 56.2
$ws.Range("F80").Value = 187.4
$ws.Range("G80").Value = 'LS DUFFY'
$ws.Range("H80").Value = 'SOLO'
$ws.Range("I80").Value = 0.1066721808657292
$ws.Range("J80").Value = 5.4
$ws.Range("K80").Value = 0.003094384707287933

# Row 81
$ws.Range("A81").Value = 84
$ws.Range("B81").Value = 1.901710291787398
$ws.Range("C81").Value = 4185.8
$ws.Range("D81").Value = 0.01349240180076153
$ws.Range("E81").Value = 29.6
$ws.Range("F81").Value = 204.4
$ws.Range("G81").Value = 'MyDogaN'
$ws.Range("H81").Value = 'DUO_SUPPORT'
$ws.Range("I81").Value = 0.09128055152368671
$ws.Range("J81").Value = 18.6
$ws.Range("K81").Value = 0.008167217339014521
